# Applies per-cell profit/price corrections produced by the scheduled
# market-data runner across the ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 15
$ws.Range("H15").Value = 496.78
$ws.Range("I15").Value = 496.78
$ws.Range("K15").Value = 1490.34
$ws.Range("M15").Value = -1321.34

# ALC row 17
$ws.Range("H17").Value = 953.25
$ws.Range("I17").Value = 498
$ws.Range("J17").Value = 1105
$ws.Range("K17").Value = 1494
$ws.Range("L17").Value = 3315
$ws.Range("M17").Value = -1326
$ws.Range("N17").Value = -3651

# ALC row 33
$ws.Range("H33").Value = 835.68964
$ws.Range("I33").Value = 841.72
$ws.Range("J33").Value = 798
$ws.Range("K33").Value = 841.72
$ws.Range("L33").Value = 798
$ws.Range("M33").Value = -612.72
$ws.Range("N33").Value = -1256

# ALC row 51
$ws.Range("H51").Value = 10105299
$ws.Range("I51").Value = 22732772
$ws.Range("J51").Value = 3319.8
$ws.Range("K51").Value = 22732772
$ws.Range("L51").Value = 3319.8
$ws.Range("M51").Value = -22732288
$ws.Range("N51").Value = -4287.8

# ALC row 62
$ws.Range("H62").Value = 2920.8462
$ws.Range("I62").Value = 2209
$ws.Range("J62").Value = 3531
$ws.Range("K62").Value = 2209
$ws.Range("L62").Value = 3531
$ws.Range("M62").Value = -1585
$ws.Range("N62").Value = -4779

# ALC row 65
$ws.Range("H65").Value = 2920.8462
$ws.Range("I65").Value = 2209
$ws.Range("J65").Value = 3531
$ws.Range("K65").Value = 11045
$ws.Range("L65").Value = 17655
$ws.Range("M65").Value = -7925
$ws.Range("N65").Value = -23895

# ALC row 98
$ws.Range("H98").Value = 2219.577
$ws.Range("I98").Value = 2337.2632
$ws.Range("J98").Value = 1900.1428
$ws.Range("K98").Value = 2337.2632
$ws.Range("L98").Value = 1900.1428
$ws.Range("M98").Value = -839.2631999999999
$ws.Range("N98").Value = -4896.1428

# ALC row 112
$ws.Range("H112").Value = 3116.7856
$ws.Range("J112").Value = 3386.25
$ws.Range("L112").Value = 10158.75
$ws.Range("N112").Value = -12374.75

# ALC row 122
$ws.Range("H122").Value = 2219.577
$ws.Range("I122").Value = 2337.2632
$ws.Range("J122").Value = 1900.1428
$ws.Range("K122").Value = 7011.7896
$ws.Range("L122").Value = 5700.428400000001
$ws.Range("M122").Value = -4561.7896
$ws.Range("N122").Value = -10600.4284

# ALC row 123
$ws.Range("H123").Value = 22962.223
$ws.Range("J123").Value = 22962.223
$ws.Range("L123").Value = 22962.223
$ws.Range("N123").Value = -32762.223

# ALC row 137
$ws.Range("H137").Value = 1576.0222
$ws.Range("I137").Value = 2205.158
$ws.Range("J137").Value = 1116.2693
$ws.Range("K137").Value = 6615.474
$ws.Range("L137").Value = 3348.8079
$ws.Range("M137").Value = -4065.474
$ws.Range("N137").Value = -8448.8079

# ALC row 138
$ws.Range("H138").Value = 2859.8208
$ws.Range("I138").Value = 1848.5358
$ws.Range("J138").Value = 3585.8718
$ws.Range("K138").Value = 5545.607400000001
$ws.Range("L138").Value = 10757.6154
$ws.Range("M138").Value = -405.6074000000008
$ws.Range("N138").Value = -21037.6154

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 13519.54
$ws.Range("I32").Value = 10474.402
$ws.Range("J32").Value = 27391.834
$ws.Range("K32").Value = 10474.402
$ws.Range("L32").Value = 27391.834
$ws.Range("M32").Value = -10187.402
$ws.Range("N32").Value = -27965.834

# ARM row 45
$ws.Range("H45").Value = 836.6842
$ws.Range("I45").Value = 824.875
$ws.Range("J45").Value = 899.6667
$ws.Range("K45").Value = 824.875
$ws.Range("L45").Value = 899.6667
$ws.Range("M45").Value = -447.875
$ws.Range("N45").Value = -1653.6667

$ws = $wb.Worksheets.Item("CRP")
# CRP row 14
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = $null

# CRP row 58
$ws.Range("H58").Value = 4075.4707
$ws.Range("I58").Value = 1998.4783
$ws.Range("J58").Value = 8418.272
$ws.Range("K58").Value = 1998.4783
$ws.Range("L58").Value = 8418.272
$ws.Range("M58").Value = -1795.4783
$ws.Range("N58").Value = -8824.272

# CRP row 132
$ws.Range("H132").Value = 2478
$ws.Range("I132").Value = 1869.5385
$ws.Range("J132").Value = 3356.889
$ws.Range("K132").Value = 5608.6155
$ws.Range("L132").Value = 10070.667
$ws.Range("M132").Value = -3078.6155
$ws.Range("N132").Value = -15130.667

# CRP row 136
$ws.Range("H136").Value = 4075.4707
$ws.Range("I136").Value = 1998.4783
$ws.Range("J136").Value = 8418.272
$ws.Range("K136").Value = 5995.4349
$ws.Range("L136").Value = 25254.816
$ws.Range("M136").Value = -3445.4349
$ws.Range("N136").Value = -30354.816

$ws = $wb.Worksheets.Item("CUL")
# CUL row 82
$ws.Range("H82").Value = 4750
$ws.Range("I82").Value = 1500
$ws.Range("K82").Value = 4500
$ws.Range("M82").Value = -4094

# CUL row 85
$ws.Range("H85").Value = 4750
$ws.Range("I85").Value = 1500
$ws.Range("K85").Value = 4500
$ws.Range("M85").Value = -3096

# CUL row 86
$ws.Range("H86").Value = 1236.75
$ws.Range("I86").Value = 482.33334
$ws.Range("J86").Value = 3500
$ws.Range("K86").Value = 1447.00002
$ws.Range("L86").Value = 10500
$ws.Range("M86").Value = -261.0000199999999
$ws.Range("N86").Value = -12872

# CUL row 89
$ws.Range("H89").Value = 1236.75
$ws.Range("I89").Value = 482.33334
$ws.Range("J89").Value = 3500
$ws.Range("K89").Value = 4341.00006
$ws.Range("L89").Value = 31500
$ws.Range("M89").Value = 1586.99994
$ws.Range("N89").Value = -43356

# CUL row 131
$ws.Range("H131").Value = 1649.1852
$ws.Range("I131").Value = 3820
$ws.Range("J131").Value = 1377.8334
$ws.Range("K131").Value = 11460
$ws.Range("L131").Value = 4133.5002
$ws.Range("M131").Value = -6420
$ws.Range("N131").Value = -14213.5002

# CUL row 132
$ws.Range("H132").Value = 1263.7826
$ws.Range("I132").Value = 939.4667
$ws.Range("J132").Value = 1871.875
$ws.Range("K132").Value = 8455.2003
$ws.Range("L132").Value = 16846.875
$ws.Range("M132").Value = -5925.2003
$ws.Range("N132").Value = -21906.875

$ws = $wb.Worksheets.Item("GSM")
# GSM row 12
$ws.Range("H12").Value = 1002
$ws.Range("I12").Value = 1003
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 1003
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = -863
$ws.Range("N12").Value = -1280

# GSM row 102
$ws.Range("H102").Value = 3670
$ws.Range("I102").Value = 4954
$ws.Range("J102").Value = 2899.6
$ws.Range("K102").Value = 4954
$ws.Range("L102").Value = 2899.6
$ws.Range("M102").Value = -3332
$ws.Range("N102").Value = -6143.6

# GSM row 132
$ws.Range("H132").Value = 3037.5625
$ws.Range("I132").Value = 2102.9092
$ws.Range("J132").Value = 3527.1428
$ws.Range("K132").Value = 6308.7276
$ws.Range("L132").Value = 10581.4284
$ws.Range("M132").Value = -3778.7276
$ws.Range("N132").Value = -15641.4284

$ws = $wb.Worksheets.Item("LTW")
# LTW row 16
$ws.Range("H16").Value = 630412.25
$ws.Range("I16").Value = 1001959.4
$ws.Range("J16").Value = 11167
$ws.Range("K16").Value = 1001959.4
$ws.Range("L16").Value = 11167
$ws.Range("M16").Value = -1001789.4
$ws.Range("N16").Value = -11507

# LTW row 18
$ws.Range("H18").Value = 5000
$ws.Range("I18").Value = 5000
$ws.Range("K18").Value = 5000
$ws.Range("M18").Value = -4828

# LTW row 40
$ws.Range("H40").Value = 4230
$ws.Range("I40").Value = 3510.7144
$ws.Range("J40").Value = 6747.5
$ws.Range("K40").Value = 3510.7144
$ws.Range("L40").Value = 6747.5
$ws.Range("M40").Value = -3374.7144
$ws.Range("N40").Value = -7019.5

# LTW row 120
$ws.Range("H120").Value = 33165.668
$ws.Range("J120").Value = 29399.5
$ws.Range("L120").Value = 29399.5
$ws.Range("N120").Value = -39075.5

# LTW row 132
$ws.Range("H132").Value = 38687.594
$ws.Range("I132").Value = 54190.25
$ws.Range("J132").Value = 12849.833
$ws.Range("K132").Value = 162570.75
$ws.Range("L132").Value = 38549.499
$ws.Range("M132").Value = -160040.75
$ws.Range("N132").Value = -43609.499

$ws = $wb.Worksheets.Item("WVR")
# WVR row 82
$ws.Range("H82").Value = 16000
$ws.Range("J82").Value = 16000
$ws.Range("L82").Value = 16000
$ws.Range("N82").Value = -16766

# WVR row 85
$ws.Range("H85").Value = 16000
$ws.Range("J85").Value = 16000
$ws.Range("L85").Value = 16000
$ws.Range("N85").Value = -18652

# WVR row 126
$ws.Range("H126").Value = 3754.077
$ws.Range("I126").Value = 3430.3
$ws.Range("J126").Value = 4833.3335
$ws.Range("K126").Value = 10290.9
$ws.Range("L126").Value = 14500.0005
$ws.Range("M126").Value = -7820.900000000001
$ws.Range("N126").Value = -19440.0005
